$wb = $excel.ActiveWorkbook

# Data for the new row 47 on each sheet, in sheet (tab) order:
# 1 ROW35-FE-LIFTER, 2 ROW35-MID-LIFTER, 3 ROW02-FE-LIFTER, 4 ROW02-MID-LIFTER
$rowsData = @(
    @{
        A = 45747.84041853009
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x76"
        E = "0xd"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 374
        I = 13
    },
    @{
        A = 45747.69308607639
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x76"
        E = "0xe"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 374
        I = 14
    },
    @{
        A = 45747.83410693287
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x76"
        E = "0x3"
        F = 400
        G = [double]"5.68631262647114e+23"
        H = 374
        I = 3
    },
    @{
        A = 45747.8890490625
        B = "0x01,0x90"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x01,0x76"
        E = "0x3"
        F = 400
        G = [double]"9.85046333984776e+23"
        H = 374
        I = 3
    }
)

for ($i = 0; $i -lt $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)
    $data = $rowsData[$i]

    $rowNum = 47

    $ws.Cells.Item($rowNum, 1).Value = $data.A
    $ws.Cells.Item($rowNum, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($rowNum, 2).Value = $data.B
    $ws.Cells.Item($rowNum, 3).Value = $data.C
    $ws.Cells.Item($rowNum, 4).Value = $data.D
    $ws.Cells.Item($rowNum, 5).Value = $data.E
    $ws.Cells.Item($rowNum, 6).Value = $data.F
    $ws.Cells.Item($rowNum, 7).Value = $data.G
    $ws.Cells.Item($rowNum, 8).Value = $data.H
    $ws.Cells.Item($rowNum, 9).Value = $data.I
}
